# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The detail table (rows 16-23, cols B:G) is re-sorted so the two workers'
# late-payment periods run in ascending order (2502, 2503, 2504, 2505),
# alternating between worker 1143366687 (ALISON ELENA MONTENEGRO LEFRANC)
# and worker 73107228 (RUBEN DARIO MONTENEGRO MERCADO).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: B=Tipo Doc, C=N Doc, D=Nombre, E=Periodo Mora, F=Valor Mora, G=Salario Basico

$rows = @(
    @{ Row = 16; Tipo = "CC"; Doc = "1143366687"; Nombre = "ALISON ELENA MONTENEGRO LEFRANC"; Periodo = "2502"; Mora = 52000;  Salario = 1300000 },
    @{ Row = 17; Tipo = "CC"; Doc = "73107228";   Nombre = "RUBEN DARIO MONTENEGRO MERCADO";   Periodo = "2502"; Mora = 140000; Salario = 3500000 },
    @{ Row = 18; Tipo = "CC"; Doc = "1143366687"; Nombre = "ALISON ELENA MONTENEGRO LEFRANC"; Periodo = "2503"; Mora = 52000;  Salario = 1300000 },
    @{ Row = 19; Tipo = "CC"; Doc = "73107228";   Nombre = "RUBEN DARIO MONTENEGRO MERCADO";   Periodo = "2503"; Mora = 140000; Salario = 3500000 },
    @{ Row = 20; Tipo = "CC"; Doc = "1143366687"; Nombre = "ALISON ELENA MONTENEGRO LEFRANC"; Periodo = "2504"; Mora = 52000;  Salario = 1300000 },
    @{ Row = 21; Tipo = "CC"; Doc = "73107228";   Nombre = "RUBEN DARIO MONTENEGRO MERCADO";   Periodo = "2504"; Mora = 140000; Salario = 3500000 },
    @{ Row = 22; Tipo = "CC"; Doc = "1143366687"; Nombre = "ALISON ELENA MONTENEGRO LEFRANC"; Periodo = "2505"; Mora = 45066;  Salario = 1300000 },
    @{ Row = 23; Tipo = "CC"; Doc = "73107228";   Nombre = "RUBEN DARIO MONTENEGRO MERCADO";   Periodo = "2505"; Mora = 121334; Salario = 3500000 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 2).Value = $r.Tipo
    $ws.Cells.Item($r.Row, 3).Value = $r.Doc
    $ws.Cells.Item($r.Row, 4).Value = $r.Nombre
    $ws.Cells.Item($r.Row, 5).Value = $r.Periodo
    $ws.Cells.Item($r.Row, 6).Value = $r.Mora
    $ws.Cells.Item($r.Row, 7).Value = $r.Salario
}
